$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column H mirrors the existing header style (bold, bordered,
# centered/top-aligned) by copying G1's formatting onto H1, then writes
# the header text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Save flag values for each data row.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 1
$ws.Range("H6").Value = 0
